$d = $word.ActiveDocument

# Locate the target sentence region precisely by searching for stable text.
$full = $d.Content.Text
$accIdx = $full.IndexOf("accounts for Walmart project, we can create an Walmart ")
if ($accIdx -lt 0) { throw "target text not found" }

$accountsStart = $accIdx
$accountsEnd   = $accountsStart + 8            # "accounts" is 8 chars
$orgIdx        = $full.IndexOf("organizational unit", $accountsEnd)
if ($orgIdx -lt 0) { throw "organizational not found" }

# --- Step 1: temporarily bold the very last character of "accounts" and the
# very first character of "organizational" so that the forthcoming text edit
# (which merges neighbouring same-formatted runs) cannot bleed into them.
$blockerBefore = $d.Range($accountsEnd - 1, $accountsEnd)
$blockerBefore.Bold = 1

$blockerAfter = $d.Range($orgIdx, $orgIdx + 1)
$blockerAfter.Bold = 1

# --- Step 2: rewrite "an Walmart" -> "a Walmart" inside the run that sits
# between "accounts" and "organizational".
$full2 = $d.Content.Text
$anIdx = $full2.IndexOf("an Walmart", $accountsEnd)
$r = $d.Range($anIdx, $anIdx + 2)
$r.Text = "a"

# --- Step 3: undo the temporary bold markers.
$full3 = $d.Content.Text
$accountsEnd2 = $full3.IndexOf("accounts") + 8
$d.Range($accountsEnd2 - 1, $accountsEnd2).Bold = 0

$orgIdx2 = $full3.IndexOf("organizational")
$d.Range($orgIdx2, $orgIdx2 + 1).Bold = 0

# --- Step 4: split "a" itself into its own run (matching the target XML
# which has "a" as a standalone <w:r>).
$full4 = $d.Content.Text
$aIdx = $full4.IndexOf("create ") + 7
$rA = $d.Range($aIdx, $aIdx + 1)
$rA.Bold = 1
$rA.Bold = 0
